$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 55.59510733333332
$ws.Range("H2").Value = 166.785322
$ws.Range("I2").Value = 0.4537221086682116
$ws.Range("J2").Value = 0.4537221086682116
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.822718
$ws.Range("N2").Value = 8.468154
$ws.Range("O2").Value = 0.0739039825466893
$ws.Range("P2").Value = 0.0739039825466893
$ws.Range("Q2").Value = 156.929310181732
$ws.Range("R2").Value = 1412.363791635588
$ws.Range("S2").Value = 0.03353187080006258
$ws.Range("T2").Value = 0.03353187080006257

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 55.59510733333332
$ws.Range("H3").Value = 166.785322
$ws.Range("I3").Value = 0.4537221086682116
$ws.Range("J3").Value = 0.4537221086682116
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.042448
$ws.Range("N3").Value = 15.127344
$ws.Range("O3").Value = 0.1320206230252502
$ws.Range("P3").Value = 0.1320206230252502
$ws.Range("Q3").Value = 280.335437782752
$ws.Range("R3").Value = 2523.018940044768
$ws.Range("S3").Value = 0.0599006754667076
$ws.Range("T3").Value = 0.05990067546670759

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 55.59510733333332
$ws.Range("H4").Value = 166.785322
$ws.Range("I4").Value = 0.4537221086682116
$ws.Range("J4").Value = 0.4537221086682116
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 28.06359933333333
$ws.Range("N4").Value = 84.190798
$ws.Range("O4").Value = 0.7347569807993387
$ws.Range("P4").Value = 0.7347569807993388
$ws.Range("Q4").Value = 1560.198817096328
$ws.Range("R4").Value = 14041.78935386695
$ws.Range("S4").Value = 0.3333754866869647
$ws.Range("T4").Value = 0.3333754866869647

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 55.59510733333332
$ws.Range("H5").Value = 166.785322
$ws.Range("I5").Value = 0.4537221086682116
$ws.Range("J5").Value = 0.4537221086682116
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.265631
$ws.Range("N5").Value = 6.796893000000001
$ws.Range("O5").Value = 0.05931841362872176
$ws.Range("P5").Value = 0.05931841362872176
$ws.Range("Q5").Value = 125.9579986227273
$ws.Range("R5").Value = 1133.621987604546
$ws.Range("S5").Value = 0.02691407571447682
$ws.Range("T5").Value = 0.02691407571447682

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.654659
$ws.Range("H6").Value = 4.963977
$ws.Range("I6").Value = 0.01350398275347337
$ws.Range("J6").Value = 0.01350398275347337
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.822718
$ws.Range("N6").Value = 8.468154
$ws.Range("O6").Value = 0.0739039825466893
$ws.Range("P6").Value = 0.0739039825466893
$ws.Range("Q6").Value = 4.670635743161999
$ws.Range("R6").Value = 42.035721688458
$ws.Range("S6").Value = 0.0009979981057234896
$ws.Range("T6").Value = 0.0009979981057234894

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.654659
$ws.Range("H7").Value = 4.963977
$ws.Range("I7").Value = 0.01350398275347337
$ws.Range("J7").Value = 0.01350398275347337
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.042448
$ws.Range("N7").Value = 15.127344
$ws.Range("O7").Value = 0.1320206230252502
$ws.Range("P7").Value = 0.1320206230252502
$ws.Range("Q7").Value = 8.343531965232
$ws.Range("R7").Value = 75.091787687088
$ws.Range("S7").Value = 0.001782804216435789
$ws.Range("T7").Value = 0.001782804216435789

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.654659
$ws.Range("H8").Value = 4.963977
$ws.Range("I8").Value = 0.01350398275347337
$ws.Range("J8").Value = 0.01350398275347337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 28.06359933333333
$ws.Range("N8").Value = 84.190798
$ws.Range("O8").Value = 0.7347569807993387
$ws.Range("P8").Value = 0.7347569807993388
$ws.Range("Q8").Value = 46.435687209294
$ws.Range("R8").Value = 417.921184883646
$ws.Range("S8").Value = 0.009922145596708436
$ws.Range("T8").Value = 0.009922145596708436

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.654659
$ws.Range("H9").Value = 4.963977
$ws.Range("I9").Value = 0.01350398275347337
$ws.Range("J9").Value = 0.01350398275347337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.265631
$ws.Range("N9").Value = 6.796893000000001
$ws.Range("O9").Value = 0.05931841362872176
$ws.Range("P9").Value = 0.05931841362872176
$ws.Range("Q9").Value = 3.748846724829
$ws.Range("R9").Value = 33.739620523461
$ws.Range("S9").Value = 0.0008010348346056586
$ws.Range("T9").Value = 0.0008010348346056585

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 62.05924166666666
$ws.Range("H10").Value = 186.177725
$ws.Range("I10").Value = 0.5064771225734745
$ws.Range("J10").Value = 0.5064771225734744
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.822718
$ws.Range("N10").Value = 8.468154
$ws.Range("O10").Value = 0.0739039825466893
$ws.Range("P10").Value = 0.0739039825466893
$ws.Range("Q10").Value = 175.17573851885
$ws.Range("R10").Value = 1576.58164666965
$ws.Range("S10").Value = 0.03743067642696748
$ws.Range("T10").Value = 0.03743067642696747

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 62.05924166666666
$ws.Range("H11").Value = 186.177725
$ws.Range("I11").Value = 0.5064771225734745
$ws.Range("J11").Value = 0.5064771225734744
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.042448
$ws.Range("N11").Value = 15.127344
$ws.Range("O11").Value = 0.1320206230252502
$ws.Range("P11").Value = 0.1320206230252502
$ws.Range("Q11").Value = 312.9304990236
$ws.Range("R11").Value = 2816.3744912124
$ws.Range("S11").Value = 0.06686542527018614
$ws.Range("T11").Value = 0.06686542527018613

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 62.05924166666666
$ws.Range("H12").Value = 186.177725
$ws.Range("I12").Value = 0.5064771225734745
$ws.Range("J12").Value = 0.5064771225734744
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 28.06359933333333
$ws.Range("N12").Value = 84.190798
$ws.Range("O12").Value = 0.7347569807993387
$ws.Range("P12").Value = 0.7347569807993388
$ws.Range("Q12").Value = 1741.605693063839
$ws.Range("R12").Value = 15674.45123757455
$ws.Range("S12").Value = 0.3721376014260228
$ws.Range("T12").Value = 0.3721376014260228

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 62.05924166666666
$ws.Range("H13").Value = 186.177725
$ws.Range("I13").Value = 0.5064771225734745
$ws.Range("J13").Value = 0.5064771225734744
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.265631
$ws.Range("N13").Value = 6.796893000000001
$ws.Range("O13").Value = 0.05931841362872176
$ws.Range("P13").Value = 0.05931841362872176
$ws.Range("Q13").Value = 140.6033417564917
$ws.Range("R13").Value = 1265.430075808425
$ws.Range("S13").Value = 0.03004341945029818
$ws.Range("T13").Value = 0.03004341945029817

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.222176333333334
$ws.Range("H14").Value = 9.666529000000001
$ws.Range("I14").Value = 0.02629678600484052
$ws.Range("J14").Value = 0.02629678600484052
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.822718
$ws.Range("N14").Value = 8.468154
$ws.Range("O14").Value = 0.0739039825466893
$ws.Range("P14").Value = 0.0739039825466893
$ws.Range("Q14").Value = 9.095295135274
$ws.Range("R14").Value = 81.857656217466
$ws.Range("S14").Value = 0.001943437213935757
$ws.Range("T14").Value = 0.001943437213935757

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.222176333333334
$ws.Range("H15").Value = 9.666529000000001
$ws.Range("I15").Value = 0.02629678600484052
$ws.Range("J15").Value = 0.02629678600484052
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.042448
$ws.Range("N15").Value = 15.127344
$ws.Range("O15").Value = 0.1320206230252502
$ws.Range("P15").Value = 0.1320206230252502
$ws.Range("Q15").Value = 16.247656607664
$ws.Range("R15").Value = 146.228909468976
$ws.Range("S15").Value = 0.003471718071920727
$ws.Range("T15").Value = 0.003471718071920726

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.222176333333334
$ws.Range("H16").Value = 9.666529000000001
$ws.Range("I16").Value = 0.02629678600484052
$ws.Range("J16").Value = 0.02629678600484052
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 28.06359933333333
$ws.Range("N16").Value = 84.190798
$ws.Range("O16").Value = 0.7347569807993387
$ws.Range("P16").Value = 0.7347569807993388
$ws.Range("Q16").Value = 90.42586560001578
$ws.Range("R16").Value = 813.832790400142
$ws.Range("S16").Value = 0.01932174708964293
$ws.Range("T16").Value = 0.01932174708964293

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.222176333333334
$ws.Range("H17").Value = 9.666529000000001
$ws.Range("I17").Value = 0.02629678600484052
$ws.Range("J17").Value = 0.02629678600484052
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.265631
$ws.Range("N17").Value = 6.796893000000001
$ws.Range("O17").Value = 0.05931841362872176
$ws.Range("P17").Value = 0.05931841362872176
$ws.Range("Q17").Value = 7.300262588266335
$ws.Range("R17").Value = 65.70236329439702
$ws.Range("S17").Value = 0.001559883629341112
$ws.Range("T17").Value = 0.001559883629341111
